# "Updated stats for Apr 3"
#
# Column I ("US Actual Confirmed Cases") holds actual reported numbers for
# past days (plain input cells, fill style 14) and forecast formulas for
# future days (fill style 15, e.g. I37 = I36*(1+AVERAGE(M34:M36))).
#
# Apr 3's actual case count (277,161) is now known, so I37 stops being a
# forecast and becomes a hard-coded actual value, like I36 before it. That
# also flips its fill color to the "actual" style. Every later forecast
# cell (I38:I49) that chains off I37 recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give I37 the same "actual value" formatting as the other real-data cells
# (e.g. I36) before overwriting its forecast formula with the real number.
$ws.Range("I36").Copy()
$ws.Range("I37").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I37").Value = 277161

# Leave the selection where Excel would after typing the value and
# pressing Enter - one row down, on I38.
$ws.Range("I38").Select()
